$p = $ppt.ActivePresentation

# The deck ships with two theme parts:
#   - theme1.xml = plain "Office Theme" palette (name "Office"), wired only
#     to the notes master's relationship.
#   - theme2.xml = the "Integral" / "Red Violet" palette, wired to the
#     slide master (and, transitively, the presentation + every slide).
#
# The authored change swaps the two parts' contents, so the slide master
# ends up using the plain Office palette while the (non-slide-facing)
# other theme becomes Red Violet.  Font scheme and format scheme (fills /
# lines / effects) are already byte-identical between the two themes, so
# the only observable difference to reproduce on the theme that actually
# drives the presentation's look is its 12-slot colour scheme (and the
# theme / colour-scheme display names, best-effort).

$sm = $p.SlideMaster
$theme = $sm.Theme

# Best-effort: keep the human-readable names in sync with the palette
# below (not all hosts persist these, but setting them is harmless).
$theme.Name = "Office Theme"
$sm.ColorScheme.Name = "Office"

$colors = $theme.ThemeColorScheme

# Standard Office theme colour scheme, as COM "RGB" long values (0xBBGGRR):
#   1  dk1      000000
#   2  lt1      FFFFFF
#   3  dk2      44546A
#   4  lt2      E7E6E6
#   5  accent1  5B9BD5
#   6  accent2  ED7D31
#   7  accent3  A5A5A5
#   8  accent4  FFC000
#   9  accent5  4472C4
#   10 accent6  70AD47
#   11 hlink    0563C1
#   12 folHlink 954F72
$colors.Item(1).RGB = 0
$colors.Item(2).RGB = 16777215
$colors.Item(3).RGB = 6968388
$colors.Item(4).RGB = 15132391
$colors.Item(5).RGB = 13998939
$colors.Item(6).RGB = 3243501
$colors.Item(7).RGB = 10855845
$colors.Item(8).RGB = 49407
$colors.Item(9).RGB = 12874308
$colors.Item(10).RGB = 4697456
$colors.Item(11).RGB = 12673797
$colors.Item(12).RGB = 7491477
